# Automatic update of files.
# Applies a data correction that:
#   - swaps the full record contents between rows 12 and 14 (two
#     mis-ordered observations got their row positions exchanged back),
#   - swaps the full record contents between rows 18 and 19 (same kind
#     of mix-up),
#   - reorders the two names listed in the "Observatörer" comment on
#     row 21.
#
# Each row is a complete observation record (Id, taxon info, coordinates,
# free-text comments, etc.). Rather than re-deriving the per-column deltas,
# we read each source row's values with Range.Value2 and write them into
# the other row (and vice versa), touching only the cells that actually
# differ so untouched cells are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($ws, $col, $r1, $r2) {
    $a1 = $col + $r1
    $a2 = $col + $r2
    $v1 = $ws.Range($a1).Value2
    $v2 = $ws.Range($a2).Value2

    if ($null -eq $v1) {
        $ws.Range($a2).ClearContents()
    } else {
        $ws.Range($a2).Value2 = $v1
    }

    if ($null -eq $v2) {
        $ws.Range($a1).ClearContents()
    } else {
        $ws.Range($a1).Value2 = $v2
    }
}

# Columns whose values differ between the two rows in each swapped pair.
$cols = @("A","B","D","E","F","G","H","M","P","Q","R","AC")

foreach ($col in $cols) {
    Swap-Cell $ws $col 12 14
}

foreach ($col in $cols) {
    Swap-Cell $ws $col 18 19
}

# Reorder the observer names on row 21.
$ws.Range("AX21").Value2 = "Anna-Lena Thommson, Lars-Erik Nilsson"

"Applied row swaps 12<->14, 18<->19 and updated AX21."
